$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @(B, C, D, E, DisNumericLooking)  ($null = leave cell unchanged)
$changes = @{
  2 = @($null, $null, "26.750.03", "  +0.33%  ", 0)
  3 = @($null, $null, "1.603.37", "  +0.32%  ", 0)
  4 = @($null, $null, $null, "  +0.15%  ", 0)
  5 = @($null, $null, "211.98", "  +0.15%  ", 1)
  6 = @($null, $null, $null, "  -0.01%  ", 0)
  7 = @($null, $null, $null, "  +0.20%  ", 0)
  9 = @($null, $null, $null, "  +0.15%  ", 0)
  10 = @($null, $null, "19.73", "  +1.03%  ", 1)
  11 = @($null, $null, "0.0848", "  +0.92%  ", 1)
  12 = @($null, $null, "1.827.10", "  +0.19%  ", 0)
  13 = @($null, $null, "1.608.92", "  +0.87%  ", 0)
  14 = @($null, $null, "4.08", "  +1.27%  ", 1)
  15 = @($null, $null, $null, "  +0.52%  ", 0)
  16 = @($null, $null, "65.09", "  +0.03%  ", 1)
  17 = @($null, $null, "0.0₃0742", "  +0.39%  ", 0)
  18 = @($null, $null, "210.61", "  +0.71%  ", 1)
  19 = @($null, $null, $null, "  +0.19%  ", 0)
  20 = @($null, $null, "7.15", "  +1.58%  ", 1)
  21 = @($null, $null, $null, "  +0.58%  ", 0)
  22 = @($null, $null, $null, "  -4.91%  ", 0)
  23 = @($null, $null, "9.08", "  +0.91%  ", 1)
  24 = @($null, $null, "143.81", "  -0.37%  ", 1)
  25 = @($null, $null, $null, "  +0.06%  ", 0)
  26 = @($null, $null, "7.11", "  +0.00%  ", 1)
  27 = @($null, $null, $null, "  -0.19%  ", 0)
  28 = @($null, $null, $null, "  +0.70%  ", 0)
  29 = @($null, $null, $null, "  -1.26%  ", 0)
  30 = @($null, $null, $null, "  -0.08%  ", 0)
  31 = @($null, $null, $null, "  +0.95%  ", 0)
  32 = @($null, $null, $null, "  +0.52%  ", 0)
  33 = @($null, $null, "1.294.95", "  +0.85%  ", 0)
  34 = @($null, $null, $null, "  +1.01%  ", 0)
  35 = @("WEMIXToken", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix", "1.21", "  +16.48%  ", 1)
  36 = @("LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "1.49", "  +0.40%  ", 1)
  37 = @($null, $null, "0.593", "  -4.39%  ", 1)
  38 = @($null, $null, $null, "  -0.13%  ", 0)
  39 = @($null, $null, "0.833", "  -0.05%  ", 1)
  40 = @($null, $null, "5.45", "  -0.45%  ", 1)
  41 = @($null, $null, $null, "  -0.30%  ", 0)
  43 = @($null, $null, "63.17", "  -0.41%  ", 1)
  44 = @($null, $null, "1.740.20", "  +0.42%  ", 0)
  45 = @("Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "6.82", "  +32.91%  ", 1)
  46 = @("Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "90.44", "  -0.73%  ", 1)
  47 = @($null, $null, "1.56", "  -0.90%  ", 1)
  48 = @("BabyDogeCoin", "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge", "0.0₆0103", "  -1.17%  ", 0)
  49 = @("Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.102", "  +0.84%  ", 1)
  50 = @("Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.0513", "  +0.83%  ", 1)
  51 = @("EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "7.58", "  +2.97%  ", 1)
}

foreach ($row in $changes.Keys) {
  $vals = $changes[$row]
  if ($vals[0] -ne $null) { $ws.Cells.Item([int]$row, 2).Value = $vals[0] }
  if ($vals[1] -ne $null) { $ws.Cells.Item([int]$row, 3).Value = $vals[1] }
  if ($vals[2] -ne $null) {
    $dCell = $ws.Cells.Item([int]$row, 4)
    if ($vals[4] -eq 1) {
      # Value reads as a number; force text storage, then restore default 'General' style
      $dCell.NumberFormat = '@'
      $dCell.Value = $vals[2]
      $dCell.NumberFormat = 'General'
      $dCell.Style = 'Normal'
    } else {
      $dCell.Value = $vals[2]
    }
  }
  if ($vals[3] -ne $null) { $ws.Cells.Item([int]$row, 5).Value = $vals[3] }
}
